$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bring in formatting for the two new rows by copying the existing
# --- similarly-bordered rows, then tweak borders to match the target layout.

# Row 6 inherits the "top+bottom thin" look of row 5 (A5:E5), then we drop
# the bottom edge so only the top separator remains.
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A6:E6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:E6").Borders.Item(9).LineStyle = -4142

# Row 7 (the new last row) inherits row 4's look (no top / thin bottom),
# then we remove the remaining bottom edge so the row ends up borderless.
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:E7").Borders.Item(9).LineStyle = -4142

$ws.Application.CutCopyMode = 0

# --- Fill in the new cell values (order matches how the shared strings
# --- table was originally appended to).
$ws.Range("C6").Value = " Yay! Good luck, [hero] and\n[partner]!"
$ws.Range("A6").Value = "SCRIPT/T01P02A/um1409.ssb"
$ws.Range("D6").Value = " Ура! Удачи вам, [hero]\nи [partner]!"
$ws.Range("E6").Value = " Ôñà! Ôäàœé âàí, [hero]\né [partner]!"
$ws.Range("A7").Value = "SCRIPT/T01P02A/um1502.ssb"

$ws.Range("B6").Value = 171

# --- Row heights.
$ws.Rows(6).RowHeight = 48
$ws.Rows(7).RowHeight = 43.2

# --- Selection / active cell, matching the saved view state.
$ws.Range("E6").Select() | Out-Null
